# Insert a new price-report row for "Lane Late / Primera" (fecha 2022-07-15,
# serial 44769) into the weekly Naranja - Vega Monumental Concepción sheet.
# Inserting at row 276 pushes the former rows 276-310 down to 277-311
# (dimension grows from A1:T310 to A1:T311); all of that data is carried
# along automatically by the row insert, only the brand-new row needs its
# values populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(276).Insert()

$ws.Range("A276").Value = 11
$ws.Range("B276").Value = "Vega Monumental Concepción"
$ws.Range("C276").Value = "Bíobío"
$ws.Range("D276").Value = 44769
$ws.Range("E276").Value = 8
$ws.Range("F276").Value = "Fruta"
$ws.Range("G276").Value = 100102
$ws.Range("H276").Value = "Cítricos"
$ws.Range("I276").Value = 100102005
$ws.Range("J276").Value = "Naranja"
$ws.Range("K276").Value = "Lane Late"
$ws.Range("L276").Value = "Primera"
$ws.Range("M276").Value = 310
$ws.Range("N276").Value = 5000
$ws.Range("O276").Value = 5500
$ws.Range("P276").Value = 5242
$ws.Range("Q276").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R276").Value = "Región de O'Higgins"
$ws.Range("S276").Value = 349
$ws.Range("T276").Value = 15
